# Actualización automática de tasas-transfi.xlsx
$wb = $excel.ActiveWorkbook

# --- Hoja1: update the "Conversión del día" summary text (cell A1) ---
$hoja1 = $wb.Worksheets.Item("Hoja1")
$cell = $hoja1.Range("A1")
$text = $cell.Value2
$text = $text -replace [regex]::Escape("✅ 1000 Bs = 4.38 = 17092.95 pesos"), "✅ 1000 Bs = 4.32 = 16756.28 pesos"
$text = $text -replace [regex]::Escape("✅ 17092.95 pesos = 4.37 = 956.16 Bs"), "✅ 16756.28 pesos = 4.29 = 944.44 Bs"
$cell.Value2 = $text

# --- tasas: update the rate cells ---
$tasas = $wb.Worksheets.Item("tasas")
$tasas.Range("N10").Value = 231.376
$tasas.Range("O10").Value = 3877
$tasas.Range("N12").Value = 3905
$tasas.Range("O12").Value = 220.1
